$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. "总计" (summary) sheet: insert a 2022-Q3 row at the top of the quarter
#    table, shifting the existing quarter rows down by one.
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item(1)

# Grow the table by one row: duplicate the last data row (row 8) into row 9
# so the new row inherits the correct cell styles (bold/boxed index column).
$summary.Range("A8:D8").Copy($summary.Range("A9:D9"))

# Shift the quarter data (columns B:D) down one row at a time, bottom-up, so
# each write reads the still-untouched value above it. Column A (the 0-based
# row index) is left untouched -- it already holds the right values.
$summary.Cells.Item(9,2).Value = $summary.Cells.Item(8,2).Value2
$summary.Cells.Item(9,3).Value = $summary.Cells.Item(8,3).Value2
$summary.Cells.Item(9,4).Value = $summary.Cells.Item(8,4).Value2

$summary.Cells.Item(8,2).Value = $summary.Cells.Item(7,2).Value2
$summary.Cells.Item(8,3).Value = $summary.Cells.Item(7,3).Value2
$summary.Cells.Item(8,4).Value = $summary.Cells.Item(7,4).Value2

$summary.Cells.Item(7,2).Value = $summary.Cells.Item(6,2).Value2
$summary.Cells.Item(7,3).Value = $summary.Cells.Item(6,3).Value2
$summary.Cells.Item(7,4).Value = $summary.Cells.Item(6,4).Value2

$summary.Cells.Item(6,2).Value = $summary.Cells.Item(5,2).Value2
$summary.Cells.Item(6,3).Value = $summary.Cells.Item(5,3).Value2
$summary.Cells.Item(6,4).Value = $summary.Cells.Item(5,4).Value2

$summary.Cells.Item(5,2).Value = $summary.Cells.Item(4,2).Value2
$summary.Cells.Item(5,3).Value = $summary.Cells.Item(4,3).Value2
$summary.Cells.Item(5,4).Value = $summary.Cells.Item(4,4).Value2

$summary.Cells.Item(4,2).Value = $summary.Cells.Item(3,2).Value2
$summary.Cells.Item(4,3).Value = $summary.Cells.Item(3,3).Value2
$summary.Cells.Item(4,4).Value = $summary.Cells.Item(3,4).Value2

$summary.Cells.Item(3,2).Value = $summary.Cells.Item(2,2).Value2
$summary.Cells.Item(3,3).Value = $summary.Cells.Item(2,3).Value2
$summary.Cells.Item(3,4).Value = $summary.Cells.Item(2,4).Value2

# Row 2 becomes the brand-new 2022-Q3 entry.
$summary.Cells.Item(2,2).Value = "2022-Q3"
$summary.Cells.Item(2,3).Value = 45
$summary.Cells.Item(2,4).Value = 3.32

# Row 9's index (column A) continues the 0-based sequence.
$summary.Cells.Item(9,1).Value = 7

# ---------------------------------------------------------------------------
# 2. Insert a brand-new "2022-Q3" worksheet right after "总计" (i.e. before
#    the current "2022-Q2" sheet), cloning the layout/formatting of the
#    "2022-Q2" sheet, then filling in the 2022-Q3 fund-holding data.
# ---------------------------------------------------------------------------
$beforeSheet = $wb.Worksheets.Item(2)
$ws = $wb.Worksheets.Add($beforeSheet)
$ws.Name = "2022-Q3"

$template = $wb.Worksheets.Item("2022-Q2")
$template.Range("A1:H57").Copy($ws.Range("A1"))

# The template sheet has 56 data rows (57 incl. header); 2022-Q3 only needs
# 45, so drop the surplus rows to match the source dimension (A1:H46).
$ws.Rows("47:57").Delete()

$ws.Cells.Item(2,1).Value = 0
$ws.Cells.Item(2,2).Value = "'" + "011230"
$ws.Cells.Item(2,3).Value = "创金合信数字经济主题股票C"
$ws.Cells.Item(2,4).Value = "'" + "15.89"
$ws.Cells.Item(2,5).Value = "'" + "92.35"
$ws.Cells.Item(2,6).Value = "'" + "3.30"
$ws.Cells.Item(2,7).Value = "'" + "0.5244"
$ws.Cells.Item(2,8).Value = 10
$ws.Cells.Item(3,1).Value = 1
$ws.Cells.Item(3,2).Value = "'" + "011229"
$ws.Cells.Item(3,3).Value = "创金合信数字经济主题股票A"
$ws.Cells.Item(3,4).Value = "'" + "15.42"
$ws.Cells.Item(3,5).Value = "'" + "92.35"
$ws.Cells.Item(3,6).Value = "'" + "3.30"
$ws.Cells.Item(3,7).Value = "'" + "0.5089"
$ws.Cells.Item(3,8).Value = 10
$ws.Cells.Item(4,1).Value = 2
$ws.Cells.Item(4,2).Value = "'" + "610002"
$ws.Cells.Item(4,3).Value = "信澳精华配置混合A"
$ws.Cells.Item(4,4).Value = "'" + "3.74"
$ws.Cells.Item(4,5).Value = "'" + "80.39"
$ws.Cells.Item(4,6).Value = "'" + "5.27"
$ws.Cells.Item(4,7).Value = "'" + "0.1971"
$ws.Cells.Item(4,8).Value = 7
$ws.Cells.Item(5,1).Value = 3
$ws.Cells.Item(5,2).Value = "'" + "009476"
$ws.Cells.Item(5,3).Value = "建信食品饮料行业股票A"
$ws.Cells.Item(5,4).Value = "'" + "3.13"
$ws.Cells.Item(5,5).Value = "'" + "88.06"
$ws.Cells.Item(5,6).Value = "'" + "5.54"
$ws.Cells.Item(5,7).Value = "'" + "0.1734"
$ws.Cells.Item(5,8).Value = 9
$ws.Cells.Item(6,1).Value = 4
$ws.Cells.Item(6,2).Value = "'" + "000551"
$ws.Cells.Item(6,3).Value = "信诚幸福消费混合"
$ws.Cells.Item(6,4).Value = "'" + "6.85"
$ws.Cells.Item(6,5).Value = "'" + "85.23"
$ws.Cells.Item(6,6).Value = "'" + "2.51"
$ws.Cells.Item(6,7).Value = "'" + "0.1719"
$ws.Cells.Item(6,8).Value = 10
$ws.Cells.Item(7,1).Value = 5
$ws.Cells.Item(7,2).Value = "'" + "160624"
$ws.Cells.Item(7,3).Value = "鹏华消费领先混合"
$ws.Cells.Item(7,4).Value = "'" + "4.52"
$ws.Cells.Item(7,5).Value = "'" + "86.35"
$ws.Cells.Item(7,6).Value = "'" + "3.56"
$ws.Cells.Item(7,7).Value = "'" + "0.1609"
$ws.Cells.Item(7,8).Value = 9
$ws.Cells.Item(8,1).Value = 6
$ws.Cells.Item(8,2).Value = "'" + "360016"
$ws.Cells.Item(8,3).Value = "光大保德信行业轮动混合"
$ws.Cells.Item(8,4).Value = "'" + "4.50"
$ws.Cells.Item(8,5).Value = "'" + "84.44"
$ws.Cells.Item(8,6).Value = "'" + "3.55"
$ws.Cells.Item(8,7).Value = "'" + "0.1598"
$ws.Cells.Item(8,8).Value = 8
$ws.Cells.Item(9,1).Value = 7
$ws.Cells.Item(9,2).Value = "'" + "012640"
$ws.Cells.Item(9,3).Value = "鹏华稳健鸿利一年持有期混合A"
$ws.Cells.Item(9,4).Value = "'" + "2.61"
$ws.Cells.Item(9,5).Value = "'" + "92.98"
$ws.Cells.Item(9,6).Value = "'" + "5.17"
$ws.Cells.Item(9,7).Value = "'" + "0.1349"
$ws.Cells.Item(9,8).Value = 5
$ws.Cells.Item(10,1).Value = 8
$ws.Cells.Item(10,2).Value = "'" + "005014"
$ws.Cells.Item(10,3).Value = "泰康景泰回报混合A"
$ws.Cells.Item(10,4).Value = "'" + "8.99"
$ws.Cells.Item(10,5).Value = "'" + "34.29"
$ws.Cells.Item(10,6).Value = "'" + "1.23"
$ws.Cells.Item(10,7).Value = "'" + "0.1106"
$ws.Cells.Item(10,8).Value = 9
$ws.Cells.Item(11,1).Value = 9
$ws.Cells.Item(11,2).Value = "'" + "002472"
$ws.Cells.Item(11,3).Value = "光大保德信先进服务业灵活配置混合A"
$ws.Cells.Item(11,4).Value = "'" + "2.51"
$ws.Cells.Item(11,5).Value = "'" + "83.99"
$ws.Cells.Item(11,6).Value = "'" + "4.40"
$ws.Cells.Item(11,7).Value = "'" + "0.1104"
$ws.Cells.Item(11,8).Value = 9
$ws.Cells.Item(12,1).Value = 10
$ws.Cells.Item(12,2).Value = "'" + "160613"
$ws.Cells.Item(12,3).Value = "鹏华盛世创新混合（LOF）"
$ws.Cells.Item(12,4).Value = "'" + "2.53"
$ws.Cells.Item(12,5).Value = "'" + "92.08"
$ws.Cells.Item(12,6).Value = "'" + "4.30"
$ws.Cells.Item(12,7).Value = "'" + "0.1088"
$ws.Cells.Item(12,8).Value = 8
$ws.Cells.Item(13,1).Value = 11
$ws.Cells.Item(13,2).Value = "'" + "012770"
$ws.Cells.Item(13,3).Value = "光大保德信创新生活混合"
$ws.Cells.Item(13,4).Value = "'" + "2.71"
$ws.Cells.Item(13,5).Value = "'" + "86.69"
$ws.Cells.Item(13,6).Value = "'" + "3.46"
$ws.Cells.Item(13,7).Value = "'" + "0.0938"
$ws.Cells.Item(13,8).Value = 9
$ws.Cells.Item(14,1).Value = 12
$ws.Cells.Item(14,2).Value = "'" + "009876"
$ws.Cells.Item(14,3).Value = "天弘甄选食品饮料股票C"
$ws.Cells.Item(14,4).Value = "'" + "1.89"
$ws.Cells.Item(14,5).Value = "'" + "84.12"
$ws.Cells.Item(14,6).Value = "'" + "4.90"
$ws.Cells.Item(14,7).Value = "'" + "0.0926"
$ws.Cells.Item(14,8).Value = 7
$ws.Cells.Item(15,1).Value = 13
$ws.Cells.Item(15,2).Value = "'" + "008166"
$ws.Cells.Item(15,3).Value = "工银消费股票A"
$ws.Cells.Item(15,4).Value = "'" + "2.83"
$ws.Cells.Item(15,5).Value = "'" + "82.76"
$ws.Cells.Item(15,6).Value = "'" + "3.27"
$ws.Cells.Item(15,7).Value = "'" + "0.0925"
$ws.Cells.Item(15,8).Value = 9
$ws.Cells.Item(16,1).Value = 14
$ws.Cells.Item(16,2).Value = "'" + "012772"
$ws.Cells.Item(16,3).Value = "信澳精华配置混合C"
$ws.Cells.Item(16,4).Value = "'" + "1.58"
$ws.Cells.Item(16,5).Value = "'" + "80.39"
$ws.Cells.Item(16,6).Value = "'" + "5.27"
$ws.Cells.Item(16,7).Value = "'" + "0.0833"
$ws.Cells.Item(16,8).Value = 7
$ws.Cells.Item(17,1).Value = 15
$ws.Cells.Item(17,2).Value = "'" + "009875"
$ws.Cells.Item(17,3).Value = "天弘甄选食品饮料股票A"
$ws.Cells.Item(17,4).Value = "'" + "1.59"
$ws.Cells.Item(17,5).Value = "'" + "84.12"
$ws.Cells.Item(17,6).Value = "'" + "4.90"
$ws.Cells.Item(17,7).Value = "'" + "0.0779"
$ws.Cells.Item(17,8).Value = 7
$ws.Cells.Item(18,1).Value = 16
$ws.Cells.Item(18,2).Value = "'" + "008134"
$ws.Cells.Item(18,3).Value = "鹏华优选价值股票"
$ws.Cells.Item(18,4).Value = "'" + "1.80"
$ws.Cells.Item(18,5).Value = "'" + "92.72"
$ws.Cells.Item(18,6).Value = "'" + "4.32"
$ws.Cells.Item(18,7).Value = "'" + "0.0778"
$ws.Cells.Item(18,8).Value = 7
$ws.Cells.Item(19,1).Value = 17
$ws.Cells.Item(19,2).Value = "'" + "001524"
$ws.Cells.Item(19,3).Value = "华泰柏瑞精选回报灵活配置混合"
$ws.Cells.Item(19,4).Value = "'" + "5.56"
$ws.Cells.Item(19,5).Value = "'" + "21.47"
$ws.Cells.Item(19,6).Value = "'" + "1.10"
$ws.Cells.Item(19,7).Value = "'" + "0.0612"
$ws.Cells.Item(19,8).Value = 4
$ws.Cells.Item(20,1).Value = 18
$ws.Cells.Item(20,2).Value = "'" + "011574"
$ws.Cells.Item(20,3).Value = "鹏华领航一年持有期混合A"
$ws.Cells.Item(20,4).Value = "'" + "1.20"
$ws.Cells.Item(20,5).Value = "'" + "92.84"
$ws.Cells.Item(20,6).Value = "'" + "4.28"
$ws.Cells.Item(20,7).Value = "'" + "0.0514"
$ws.Cells.Item(20,8).Value = 7
$ws.Cells.Item(21,1).Value = 19
$ws.Cells.Item(21,2).Value = "'" + "011431"
$ws.Cells.Item(21,3).Value = "泰达宏利消费服务混合A"
$ws.Cells.Item(21,4).Value = "'" + "1.47"
$ws.Cells.Item(21,5).Value = "'" + "86.98"
$ws.Cells.Item(21,6).Value = "'" + "2.92"
$ws.Cells.Item(21,7).Value = "'" + "0.0429"
$ws.Cells.Item(21,8).Value = 7
$ws.Cells.Item(22,1).Value = 20
$ws.Cells.Item(22,2).Value = "'" + "011575"
$ws.Cells.Item(22,3).Value = "鹏华领航一年持有期混合C"
$ws.Cells.Item(22,4).Value = "'" + "0.91"
$ws.Cells.Item(22,5).Value = "'" + "92.84"
$ws.Cells.Item(22,6).Value = "'" + "4.28"
$ws.Cells.Item(22,7).Value = "'" + "0.0389"
$ws.Cells.Item(22,8).Value = 7
$ws.Cells.Item(23,1).Value = 21
$ws.Cells.Item(23,2).Value = "'" + "014864"
$ws.Cells.Item(23,3).Value = "建信食品饮料行业股票C"
$ws.Cells.Item(23,4).Value = "'" + "0.62"
$ws.Cells.Item(23,5).Value = "'" + "88.06"
$ws.Cells.Item(23,6).Value = "'" + "5.54"
$ws.Cells.Item(23,7).Value = "'" + "0.0343"
$ws.Cells.Item(23,8).Value = 9
$ws.Cells.Item(24,1).Value = 22
$ws.Cells.Item(24,2).Value = "'" + "003105"
$ws.Cells.Item(24,3).Value = "光大保德信永鑫灵活配置混合A"
$ws.Cells.Item(24,4).Value = "'" + "2.19"
$ws.Cells.Item(24,5).Value = "'" + "25.70"
$ws.Cells.Item(24,6).Value = "'" + "1.35"
$ws.Cells.Item(24,7).Value = "'" + "0.0296"
$ws.Cells.Item(24,8).Value = 7
$ws.Cells.Item(25,1).Value = 23
$ws.Cells.Item(25,2).Value = "'" + "008167"
$ws.Cells.Item(25,3).Value = "工银消费股票C"
$ws.Cells.Item(25,4).Value = "'" + "0.70"
$ws.Cells.Item(25,5).Value = "'" + "82.76"
$ws.Cells.Item(25,6).Value = "'" + "3.27"
$ws.Cells.Item(25,7).Value = "'" + "0.0229"
$ws.Cells.Item(25,8).Value = 9
$ws.Cells.Item(26,1).Value = 24
$ws.Cells.Item(26,2).Value = "'" + "009899"
$ws.Cells.Item(26,3).Value = "上银内需增长股票A"
$ws.Cells.Item(26,4).Value = "'" + "0.57"
$ws.Cells.Item(26,5).Value = "'" + "90.70"
$ws.Cells.Item(26,6).Value = "'" + "3.15"
$ws.Cells.Item(26,7).Value = "'" + "0.0180"
$ws.Cells.Item(26,8).Value = 9
$ws.Cells.Item(27,1).Value = 25
$ws.Cells.Item(27,2).Value = "'" + "001464"
$ws.Cells.Item(27,3).Value = "光大保德信鼎鑫灵活配置混合A"
$ws.Cells.Item(27,4).Value = "'" + "1.49"
$ws.Cells.Item(27,5).Value = "'" + "24.64"
$ws.Cells.Item(27,6).Value = "'" + "1.00"
$ws.Cells.Item(27,7).Value = "'" + "0.0149"
$ws.Cells.Item(27,8).Value = 10
$ws.Cells.Item(28,1).Value = 26
$ws.Cells.Item(28,2).Value = "'" + "011432"
$ws.Cells.Item(28,3).Value = "泰达宏利消费服务混合C"
$ws.Cells.Item(28,4).Value = "'" + "0.44"
$ws.Cells.Item(28,5).Value = "'" + "86.98"
$ws.Cells.Item(28,6).Value = "'" + "2.92"
$ws.Cells.Item(28,7).Value = "'" + "0.0128"
$ws.Cells.Item(28,8).Value = 7
$ws.Cells.Item(29,1).Value = 27
$ws.Cells.Item(29,2).Value = "'" + "010428"
$ws.Cells.Item(29,3).Value = "兴银策略智选混合C"
$ws.Cells.Item(29,4).Value = "'" + "0.40"
$ws.Cells.Item(29,5).Value = "'" + "90.21"
$ws.Cells.Item(29,6).Value = "'" + "3.19"
$ws.Cells.Item(29,7).Value = "'" + "0.0128"
$ws.Cells.Item(29,8).Value = 6
$ws.Cells.Item(30,1).Value = 28
$ws.Cells.Item(30,2).Value = "'" + "004189"
$ws.Cells.Item(30,3).Value = "华商消费行业股票"
$ws.Cells.Item(30,4).Value = "'" + "0.32"
$ws.Cells.Item(30,5).Value = "'" + "80.75"
$ws.Cells.Item(30,6).Value = "'" + "3.94"
$ws.Cells.Item(30,7).Value = "'" + "0.0126"
$ws.Cells.Item(30,8).Value = 4
$ws.Cells.Item(31,1).Value = 29
$ws.Cells.Item(31,2).Value = "'" + "009527"
$ws.Cells.Item(31,3).Value = "浙商汇金新兴消费灵活配置混合"
$ws.Cells.Item(31,4).Value = "'" + "0.51"
$ws.Cells.Item(31,5).Value = "'" + "39.77"
$ws.Cells.Item(31,6).Value = "'" + "2.06"
$ws.Cells.Item(31,7).Value = "'" + "0.0105"
$ws.Cells.Item(31,8).Value = 10
$ws.Cells.Item(32,1).Value = 30
$ws.Cells.Item(32,2).Value = "'" + "007308"
$ws.Cells.Item(32,3).Value = "华宝消费升级混合"
$ws.Cells.Item(32,4).Value = "'" + "0.53"
$ws.Cells.Item(32,5).Value = "'" + "85.14"
$ws.Cells.Item(32,6).Value = "'" + "1.80"
$ws.Cells.Item(32,7).Value = "'" + "0.0095"
$ws.Cells.Item(32,8).Value = 10
$ws.Cells.Item(33,1).Value = 31
$ws.Cells.Item(33,2).Value = "'" + "010711"
$ws.Cells.Item(33,3).Value = "华富国潮优选混合"
$ws.Cells.Item(33,4).Value = "'" + "0.22"
$ws.Cells.Item(33,5).Value = "'" + "89.37"
$ws.Cells.Item(33,6).Value = "'" + "3.80"
$ws.Cells.Item(33,7).Value = "'" + "0.0084"
$ws.Cells.Item(33,8).Value = 9
$ws.Cells.Item(34,1).Value = 32
$ws.Cells.Item(34,2).Value = "'" + "010427"
$ws.Cells.Item(34,3).Value = "兴银策略智选混合A"
$ws.Cells.Item(34,4).Value = "'" + "0.26"
$ws.Cells.Item(34,5).Value = "'" + "90.21"
$ws.Cells.Item(34,6).Value = "'" + "3.19"
$ws.Cells.Item(34,7).Value = "'" + "0.0083"
$ws.Cells.Item(34,8).Value = 6
$ws.Cells.Item(35,1).Value = 33
$ws.Cells.Item(35,2).Value = "'" + "015889"
$ws.Cells.Item(35,3).Value = "富国中证500基本面精选股票A"
$ws.Cells.Item(35,4).Value = "'" + "0.97"
$ws.Cells.Item(35,5).Value = "'" + "40.38"
$ws.Cells.Item(35,6).Value = "'" + "0.80"
$ws.Cells.Item(35,7).Value = "'" + "0.0078"
$ws.Cells.Item(35,8).Value = 5
$ws.Cells.Item(36,1).Value = 34
$ws.Cells.Item(36,2).Value = "'" + "013350"
$ws.Cells.Item(36,3).Value = "光大保德信先进服务业灵活配置混合C"
$ws.Cells.Item(36,4).Value = "'" + "0.17"
$ws.Cells.Item(36,5).Value = "'" + "83.99"
$ws.Cells.Item(36,6).Value = "'" + "4.40"
$ws.Cells.Item(36,7).Value = "'" + "0.0075"
$ws.Cells.Item(36,8).Value = 9
$ws.Cells.Item(37,1).Value = 35
$ws.Cells.Item(37,2).Value = "'" + "003242"
$ws.Cells.Item(37,3).Value = "创金合信量化发现灵活配置混合C"
$ws.Cells.Item(37,4).Value = "'" + "0.40"
$ws.Cells.Item(37,5).Value = "'" + "92.08"
$ws.Cells.Item(37,6).Value = "'" + "1.55"
$ws.Cells.Item(37,7).Value = "'" + "0.0062"
$ws.Cells.Item(37,8).Value = 6
$ws.Cells.Item(38,1).Value = 36
$ws.Cells.Item(38,2).Value = "'" + "003106"
$ws.Cells.Item(38,3).Value = "光大保德信永鑫灵活配置混合C"
$ws.Cells.Item(38,4).Value = "'" + "0.39"
$ws.Cells.Item(38,5).Value = "'" + "25.70"
$ws.Cells.Item(38,6).Value = "'" + "1.35"
$ws.Cells.Item(38,7).Value = "'" + "0.0053"
$ws.Cells.Item(38,8).Value = 7
$ws.Cells.Item(39,1).Value = 37
$ws.Cells.Item(39,2).Value = "'" + "012641"
$ws.Cells.Item(39,3).Value = "鹏华稳健鸿利一年持有期混合C"
$ws.Cells.Item(39,4).Value = "'" + "0.10"
$ws.Cells.Item(39,5).Value = "'" + "92.98"
$ws.Cells.Item(39,6).Value = "'" + "5.17"
$ws.Cells.Item(39,7).Value = "'" + "0.0052"
$ws.Cells.Item(39,8).Value = 5
$ws.Cells.Item(40,1).Value = 38
$ws.Cells.Item(40,2).Value = "'" + "003241"
$ws.Cells.Item(40,3).Value = "创金合信量化发现灵活配置混合A"
$ws.Cells.Item(40,4).Value = "'" + "0.32"
$ws.Cells.Item(40,5).Value = "'" + "92.08"
$ws.Cells.Item(40,6).Value = "'" + "1.55"
$ws.Cells.Item(40,7).Value = "'" + "0.0050"
$ws.Cells.Item(40,8).Value = 6
$ws.Cells.Item(41,1).Value = 39
$ws.Cells.Item(41,2).Value = "'" + "005015"
$ws.Cells.Item(41,3).Value = "泰康景泰回报混合C"
$ws.Cells.Item(41,4).Value = "'" + "0.39"
$ws.Cells.Item(41,5).Value = "'" + "34.29"
$ws.Cells.Item(41,6).Value = "'" + "1.23"
$ws.Cells.Item(41,7).Value = "'" + "0.0048"
$ws.Cells.Item(41,8).Value = 9
$ws.Cells.Item(42,1).Value = 40
$ws.Cells.Item(42,2).Value = "'" + "001730"
$ws.Cells.Item(42,3).Value = "兴银大健康灵活配置混合"
$ws.Cells.Item(42,4).Value = "'" + "0.15"
$ws.Cells.Item(42,5).Value = "'" + "90.36"
$ws.Cells.Item(42,6).Value = "'" + "3.12"
$ws.Cells.Item(42,7).Value = "'" + "0.0047"
$ws.Cells.Item(42,8).Value = 9
$ws.Cells.Item(43,1).Value = 41
$ws.Cells.Item(43,2).Value = "'" + "004456"
$ws.Cells.Item(43,3).Value = "兴银消费新趋势灵活配置混合"
$ws.Cells.Item(43,4).Value = "'" + "0.10"
$ws.Cells.Item(43,5).Value = "'" + "90.96"
$ws.Cells.Item(43,6).Value = "'" + "4.06"
$ws.Cells.Item(43,7).Value = "'" + "0.0041"
$ws.Cells.Item(43,8).Value = 8
$ws.Cells.Item(44,1).Value = 42
$ws.Cells.Item(44,2).Value = "'" + "001823"
$ws.Cells.Item(44,3).Value = "光大保德信鼎鑫灵活配置混合C"
$ws.Cells.Item(44,4).Value = "'" + "0.26"
$ws.Cells.Item(44,5).Value = "'" + "24.64"
$ws.Cells.Item(44,6).Value = "'" + "1.00"
$ws.Cells.Item(44,7).Value = "'" + "0.0026"
$ws.Cells.Item(44,8).Value = 10
$ws.Cells.Item(45,1).Value = 43
$ws.Cells.Item(45,2).Value = "'" + "015890"
$ws.Cells.Item(45,3).Value = "富国中证500基本面精选股票C"
$ws.Cells.Item(45,4).Value = "'" + "0.06"
$ws.Cells.Item(45,5).Value = "'" + "40.38"
$ws.Cells.Item(45,6).Value = "'" + "0.80"
$ws.Cells.Item(45,7).Value = "'" + "0.0005"
$ws.Cells.Item(45,8).Value = 5
$ws.Cells.Item(46,1).Value = 44
$ws.Cells.Item(46,2).Value = "'" + "015754"
$ws.Cells.Item(46,3).Value = "上银内需增长股票C"
$ws.Cells.Item(46,4).Value = "'" + "0.00"
$ws.Cells.Item(46,5).Value = "'" + "90.70"
$ws.Cells.Item(46,6).Value = "'" + "3.15"
$ws.Cells.Item(46,7).Value = 0
$ws.Cells.Item(46,8).Value = 9

# The B,D,E,F,G columns above were written with a leading apostrophe to force
# text storage (otherwise Excel reinterprets the numeric-looking strings as
# numbers/dates). Clear the resulting "Text"-format/quote-prefix styling so
# the cells fall back to the sheet's plain (unstyled) look, matching every
# other data cell.
$ws.Range("B2:G46").ClearFormats()

Write-Host "done"
